{"js": "// CIV-5793 fix bugs related to generating document for part admit states paid\n//\n// The merge-field condition that decides whether the \"How the defendant\n// will pay\" section is shown currently reads:\n//   <<cs_{responseType == 'FULL_ADMISSION' || (responseType == 'PART_ADMISSION' && howToPay != null)}>>\n// It needs an extra clause so the section is suppressed when the\n// defendant's reason for rejecting is that the claim was already paid:\n//   <<cs_{responseType == 'FULL_ADMISSION' || (responseType == 'PART_ADMISSION' && howToPay != null && whyReject!= 'ALREADY_PAID')}>>\n//\n// Locate the run containing \" != null)\" right after \"howToPay\" and extend\n// it in place so the visible text gains \" && whyReject!= \\u2018ALREADY_PAID\\u2019\"\n// before the closing parenthesis.\n\nconst body = context.document.body;\nconst matches = body.search(\"!= null)\", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items,text\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find the ' != null)' merge-field condition to update.\");\n}\n\nconst target = matches.items[0];\nconst replacement =\n  \"!= null && whyReject!= \\u2018ALREADY_PAID\\u2019)\";\ntarget.insertText(replacement, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# CIV-5793 fix bugs related to generating document for part admit states paid\n#\n# The merge-field condition that decides whether the \"How the defendant\n# will pay\" section is shown currently reads:\n#   <<cs_{responseType == 'FULL_ADMISSION' || (responseType == 'PART_ADMISSION' && howToPay != null)}>>\n# It needs an extra clause so the section is suppressed when the\n# defendant's reason for rejecting is that the claim was already paid:\n#   <<cs_{responseType == 'FULL_ADMISSION' || (responseType == 'PART_ADMISSION' && howToPay != null && whyReject!= 'ALREADY_PAID')}>>\n#\n# Build the replacement text using explicit char codes for the curly\n# quotes (typed smart quotes get auto-folded to straight ASCII quotes by\n# this shell, so [char] is used to keep the real Unicode punctuation used\n# throughout the rest of the template).\n$lq = [char]0x2018\n$rq = [char]0x2019\n\n$d = $word.ActiveDocument\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"!= null)\"\n$find.MatchCase = $true\n$find.Replacement.Text = \"!= null && whyReject!= \" + $lq + \"ALREADY_PAID\" + $rq + \")\"\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw \"Could not find the ' != null)' merge-field condition to update.\"\n}\n"}
